$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$sortRange = $ws.Range("A1:N9")
$keyRange = $ws.Range("B1:B9")

$ws.Sort.SortFields.Clear()
[void]$ws.Sort.SortFields.Add($keyRange, 0, 1, 0, 0)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 1
$ws.Sort.Apply()

$ws.Range("A1:N1").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$N`$1"
    }
}
